$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 4444.1333
$ws.Range("I64").Value = 4337.375
$ws.Range("J64").Value = 4566.143
$ws.Range("K64").Value = 4337.375
$ws.Range("L64").Value = 4566.143
$ws.Range("M64").Value = -4089.375
$ws.Range("N64").Value = -5062.143
$ws.Range("H67").Value = 4444.1333
$ws.Range("I67").Value = 4337.375
$ws.Range("J67").Value = 4566.143
$ws.Range("K67").Value = 4337.375
$ws.Range("L67").Value = 4566.143
$ws.Range("M67").Value = -3479.375
$ws.Range("N67").Value = -6282.143
$ws.Range("H74").Value = 8199.75
$ws.Range("I74").Value = 9719.6
$ws.Range("J74").Value = 5666.6665
$ws.Range("K74").Value = 9719.6
$ws.Range("L74").Value = 5666.6665
$ws.Range("M74").Value = -8783.6
$ws.Range("N74").Value = -7538.6665
$ws.Range("H76").Value = 4286.1816
$ws.Range("I76").Value = 3549.3333
$ws.Range("J76").Value = 4562.5
$ws.Range("K76").Value = 3549.3333
$ws.Range("L76").Value = 4562.5
$ws.Range("M76").Value = -3234.3333
$ws.Range("N76").Value = -5192.5
$ws.Range("H77").Value = 8199.75
$ws.Range("I77").Value = 9719.6
$ws.Range("J77").Value = 5666.6665
$ws.Range("K77").Value = 48598
$ws.Range("L77").Value = 28333.3325
$ws.Range("M77").Value = -43918
$ws.Range("N77").Value = -37693.3325
$ws.Range("H79").Value = 4286.1816
$ws.Range("I79").Value = 3549.3333
$ws.Range("J79").Value = 4562.5
$ws.Range("K79").Value = 3549.3333
$ws.Range("L79").Value = 4562.5
$ws.Range("M79").Value = -2457.3333
$ws.Range("N79").Value = -6746.5
$ws.Range("H135").Value = 2769.8474
$ws.Range("I135").Value = 3201.976
$ws.Range("J135").Value = 1702.2354
$ws.Range("K135").Value = 28817.784
$ws.Range("L135").Value = 15320.1186
$ws.Range("M135").Value = -26282.784
$ws.Range("N135").Value = -20390.1186
$ws.Range("H137").Value = 196387.89
$ws.Range("I137").Value = 411371.8
$ws.Range("J137").Value = 1558.6875
$ws.Range("K137").Value = 1234115.4
$ws.Range("L137").Value = 4676.0625
$ws.Range("M137").Value = -1231565.4
$ws.Range("N137").Value = -9776.0625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3644.5
$ws.Range("I32").Value = 3142.1292
$ws.Range("K32").Value = 3142.1292
$ws.Range("M32").Value = -2855.1292
$ws.Range("H74").Value = 5441.636
$ws.Range("I74").Value = 428.9
$ws.Range("J74").Value = 9618.916999999999
$ws.Range("K74").Value = 428.9
$ws.Range("L74").Value = 9618.916999999999
$ws.Range("M74").Value = 445.1
$ws.Range("N74").Value = -11366.917
$ws.Range("H77").Value = 5441.636
$ws.Range("I77").Value = 428.9
$ws.Range("J77").Value = 9618.916999999999
$ws.Range("K77").Value = 2144.5
$ws.Range("L77").Value = 48094.585
$ws.Range("M77").Value = 2223.5
$ws.Range("N77").Value = -56830.585

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 13910870
$ws.Range("I134").Value = 19638168
$ws.Range("J134").Value = 1716.2858
$ws.Range("K134").Value = 58914504
$ws.Range("L134").Value = 5148.857400000001
$ws.Range("M134").Value = -58911969
$ws.Range("N134").Value = -10218.8574

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6352.7437
$ws.Range("I31").Value = 927.8679
$ws.Range("J31").Value = 17853.48
$ws.Range("K31").Value = 927.8679
$ws.Range("L31").Value = 17853.48
$ws.Range("M31").Value = -632.8679
$ws.Range("N31").Value = -18443.48
$ws.Range("H34").Value = 6352.7437
$ws.Range("I34").Value = 927.8679
$ws.Range("J34").Value = 17853.48
$ws.Range("K34").Value = 927.8679
$ws.Range("L34").Value = 17853.48
$ws.Range("M34").Value = -725.8679
$ws.Range("N34").Value = -18257.48
$ws.Range("H132").Value = 9528435
$ws.Range("I132").Value = 20834744
$ws.Range("J132").Value = 7332.2104
$ws.Range("K132").Value = 62504232
$ws.Range("L132").Value = 21996.6312
$ws.Range("M132").Value = -62501702
$ws.Range("N132").Value = -27056.6312
$ws.Range("H134").Value = 5788040
$ws.Range("I134").Value = 6579807.5
$ws.Range("J134").Value = 3907591.8
$ws.Range("K134").Value = 19739422.5
$ws.Range("L134").Value = 11722775.4
$ws.Range("M134").Value = -19736887.5
$ws.Range("N134").Value = -11727845.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 855
$ws.Range("I122").Value = 645.25
$ws.Range("J122").Value = 1214.5714
$ws.Range("K122").Value = 5807.25
$ws.Range("L122").Value = 10931.1426
$ws.Range("M122").Value = -3357.25
$ws.Range("N122").Value = -15831.1426

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("M27").ClearContents()
$ws.Range("H80").Value = 5333.3335
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 5333.3335
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 5333.3335
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -7329.3335
$ws.Range("H83").Value = 5333.3335
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 5333.3335
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 26666.6675
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -36650.6675

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H115").Value = 9000
$ws.Range("I115").Value = 9000
$ws.Range("K115").Value = 9000
$ws.Range("M115").Value = -7825
$ws.Range("H136").Value = 4531.625
$ws.Range("I136").Value = 5163.2256
$ws.Range("J136").Value = 2356.111
$ws.Range("K136").Value = 15489.6768
$ws.Range("L136").Value = 7068.333
$ws.Range("M136").Value = -12939.6768
$ws.Range("N136").Value = -12168.333

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 40812076
$ws.Range("I132").Value = 37501300
$ws.Range("J132").Value = 51406570
$ws.Range("K132").Value = 112503900
$ws.Range("L132").Value = 154219710
$ws.Range("M132").Value = -112501370
$ws.Range("N132").Value = -154224770
$ws.Range("H136").Value = 13824953
$ws.Range("I136").Value = 8116139.5
$ws.Range("J136").Value = 27779832
$ws.Range("K136").Value = 24348418.5
$ws.Range("L136").Value = 83339496
$ws.Range("M136").Value = -24345868.5
$ws.Range("N136").Value = -83344596
